$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.864.78"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "1.720.62"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +1.05%  "
$ws.Range("D5").Value = "'312.81"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("D7").Value = "'0.3792"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.3505"
$ws.Range("E8").Value = "  +2.01%  "
$ws.Range("D9").Value = "'49.63"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "'1.192"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "'0.07504"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "'1.013"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "'6.371"
$ws.Range("E13").Value = "  +4.76%  "
$ws.Range("D14").Value = "'20.87"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "'6.998"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").Value = "1.730.62"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").Value = "'0.00001128"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "'0.06708"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'84.49"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("D21").Value = "'17.34"
$ws.Range("E21").Value = "  +5.38%  "
$ws.Range("D22").Value = "'6.390"
$ws.Range("E22").Value = "  +4.35%  "
$ws.Range("D23").Value = "'13.11"
$ws.Range("E23").Value = "  +9.18%  "
$ws.Range("D24").Value = "24.877.32"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").Value = "'2.452"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").Value = "'2.801"
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").Value = "'20.48"
$ws.Range("E27").Value = "  +4.72%  "
$ws.Range("D28").Value = "'151.64"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").Value = "'131.65"
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("D31").Value = "'1.174"
$ws.Range("E31").Value = "  +18.62%  "
$ws.Range("D32").Value = "'6.850"
$ws.Range("E32").Value = "  +6.33%  "
$ws.Range("D33").Value = "'4.251"
$ws.Range("E33").Value = "  +5.36%  "
$ws.Range("D34").Value = "'1.803"
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").Value = "'13.75"
$ws.Range("E35").Value = "  +10.47%  "
$ws.Range("D36").Value = "'0.08719"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("D37").Value = "'5.597"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("D38").Value = "'0.02445"
$ws.Range("E38").Value = "  +4.58%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06565"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'9.102"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "'0.2204"
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("D42").Value = "'1.274"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").Value = "'0.6461"
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "'14.00"
$ws.Range("E45").Value = "  +5.96%  "
$ws.Range("D46").Value = "'0.6175"
$ws.Range("E46").Value = "  +3.64%  "
$ws.Range("D47").Value = "'3.853"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "'2.153"
$ws.Range("E48").Value = "  +6.74%  "
$ws.Range("D49").Value = "'128.85"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "'0.07276"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").Value = "'79.89"
$ws.Range("E51").Value = "  +4.43%  "
